$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Enter values in the order that produces the shared-string table
# [f, s, d, dd] (first-use order determines the uniqueCount index).
$ws.Range("C5").Value = "f"
$ws.Range("G6").Value = "f"

$ws.Range("D9").Value = "s"

$ws.Range("N4").Value = "d"
$ws.Range("H5").Value = "d"
$ws.Range("H8").Value = "d"

$ws.Range("K3").Value = "dd"

# Restore the selected cell shown in the saved workbook.
$ws.Range("E5").Select()
